# CCS Calculations.xlsx update
#
# The "Electricity" column in the Calculations sheet was being sourced from
# Table 2's IGCC (2012) figures; it should instead be sourced from Table 2's
# PC (2012) figures. This changes:
#   - the column header label (row 13 / row 21) from "Electricity (IGCC)"
#     to "Electricity (PC)" (matching the label already used in row 2)
#   - the Capital Cost input (row 14) to pull from column E instead of C
#   - the O&M input (row 22) to pull from column E instead of C
# Downstream formulas (B17, B25) and the dependent summary sheets
# (CC-TOMCpTS, CC-EUpTCS) recalculate automatically.

$wb = $excel.ActiveWorkbook

$calc = $wb.Worksheets.Item("Calculations")

# Relabel the "Electricity" column headers to match the PC source data.
$calc.Range("B13").Value = "Electricity (PC)"
$calc.Range("B21").Value = "Electricity (PC)"

# Re-point the source formulas from the IGCC column (C) to the PC column (E).
$calc.Range("B14").Formula = "='Table 2'!E10"
$calc.Range("B22").Formula = "='Table 2'!E11"

# View-state: CC-TOMCpTS selection moves to B39:B40.
$tomcpts = $wb.Worksheets.Item("CC-TOMCpTS")
$tomcpts.Range("B39:B40").Select()

# View-state: Calculations becomes the active/front sheet with B14:B25 selected.
$calc.Activate()
$calc.Range("B14:B25").Select()
